# Auto-generated Excel COM-interop script
# Applies numeric updates to the H-N (price/profit) columns
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets,
# matching the authoritative diff of the workbook's OOXML.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 12072.154
$ws.Range("I18").Value = 12585.363
$ws.Range("K18").Value = 12585.363
$ws.Range("M18").Value = -12301.363

$ws.Range("H55").Value = 269.33334
$ws.Range("I55").Value = 134.8
$ws.Range("J55").Value = 437.5
$ws.Range("K55").Value = 134.8
$ws.Range("L55").Value = 437.5
$ws.Range("M55").Value = 79.19999999999999
$ws.Range("N55").Value = -865.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3448.6155
$ws.Range("J32").Value = 9499
$ws.Range("L32").Value = 9499
$ws.Range("N32").Value = -10073

$ws.Range("H43").Value = 45671
$ws.Range("I43").Value = 45671
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 45671
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -45358
$ws.Range("N43").Value = $null

$ws.Range("H45").Value = 3542.1667
$ws.Range("I45").Value = 2563.25
$ws.Range("K45").Value = 2563.25
$ws.Range("M45").Value = -2186.25

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null

$ws.Range("H132").Value = 2220.5
$ws.Range("I132").Value = 2220.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6661.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4131.5
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 323.75
$ws.Range("I22").Value = 323.75
$ws.Range("K22").Value = 323.75
$ws.Range("M22").Value = -150.75

$ws.Range("H34").Value = 8000
$ws.Range("I34").Value = 8000
$ws.Range("K34").Value = 8000
$ws.Range("M34").Value = -7886

$ws.Range("H94").Value = 2414
$ws.Range("I94").Value = 2414
$ws.Range("K94").Value = 2414
$ws.Range("M94").Value = -1963

$ws.Range("H134").Value = 6990.154
$ws.Range("I134").Value = 3985
$ws.Range("K134").Value = 11955
$ws.Range("M134").Value = -9420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1679999.4
$ws.Range("I3").Value = 2509999.5
$ws.Range("J3").Value = 19999
$ws.Range("K3").Value = 2509999.5
$ws.Range("L3").Value = 19999
$ws.Range("M3").Value = -2509886.5
$ws.Range("N3").Value = -20225

$ws.Range("H51").Value = 22857.143
$ws.Range("J51").Value = 23333.334
$ws.Range("L51").Value = 23333.334
$ws.Range("N51").Value = -24805.334

$ws.Range("H61").Value = 22857.143
$ws.Range("J61").Value = 23333.334
$ws.Range("L61").Value = 23333.334
$ws.Range("N61").Value = -24029.334

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2039.45
$ws.Range("J55").Value = 2700
$ws.Range("L55").Value = 8100
$ws.Range("N55").Value = -8454

$ws.Range("H111").Value = 1948.5
$ws.Range("I111").Value = 1948.5
$ws.Range("K111").Value = 5845.5
$ws.Range("M111").Value = -2778.5

$ws.Range("H131").Value = 3713.0625
$ws.Range("I131").Value = 1632
$ws.Range("J131").Value = 4659
$ws.Range("K131").Value = 4896
$ws.Range("L131").Value = 13977
$ws.Range("M131").Value = 144
$ws.Range("N131").Value = -24057

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 25004000
$ws.Range("I7").Value = 50000000
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 50000000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -49999888
$ws.Range("N7").Value = -8224

$ws.Range("H8").Value = 25004000
$ws.Range("I8").Value = 50000000
$ws.Range("J8").Value = 8000
$ws.Range("K8").Value = 50000000
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = -49999861
$ws.Range("N8").Value = -8278

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 20005
$ws.Range("I3").Value = 20005
$ws.Range("K3").Value = 20005
$ws.Range("M3").Value = -19893

$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null

$ws.Range("H15").Value = 20005
$ws.Range("I15").Value = 20005
$ws.Range("K15").Value = 20005
$ws.Range("M15").Value = -19835

$ws.Range("H21").Value = 8594
$ws.Range("I21").Value = 8125.3335
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 8125.3335
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -7951.3335
$ws.Range("N21").Value = -10348

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").Value = $null

$ws.Range("H50").Value = 17500
$ws.Range("I50").Value = 17500
$ws.Range("K50").Value = 17500
$ws.Range("M50").Value = -16863

$ws.Range("H56").Value = 37515000
$ws.Range("I56").Value = 37515000
$ws.Range("K56").Value = 37515000
$ws.Range("M56").Value = -37514309

$ws.Range("H116").Value = 100000
$ws.Range("J116").Value = 100000
$ws.Range("L116").Value = 100000
$ws.Range("N116").Value = -109178

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 1671701
$ws.Range("I3").Value = 1671701
$ws.Range("K3").Value = 1671701
$ws.Range("M3").Value = -1671587

$ws.Range("H11").Value = 9999
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 9999
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 9999
$ws.Range("M11").Value = $null
$ws.Range("N11").Value = -10283

$ws.Range("H14").Value = 3000
$ws.Range("I14").Value = 3000
$ws.Range("K14").Value = 3000
$ws.Range("M14").Value = -2832

$ws.Range("H15").Value = 10000
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").Value = $null

$ws.Range("H18").Value = 7977.769
$ws.Range("I18").Value = 25000
$ws.Range("J18").Value = 2871.1
$ws.Range("K18").Value = 25000
$ws.Range("L18").Value = 2871.1
$ws.Range("M18").Value = -24827
$ws.Range("N18").Value = -3217.1

$ws.Range("H19").Value = 6400
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 6400
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 6400
$ws.Range("M19").Value = $null
$ws.Range("N19").Value = -6748

$ws.Range("H20").Value = 650
$ws.Range("I20").Value = 650
$ws.Range("K20").Value = 650
$ws.Range("M20").Value = -410

$ws.Range("H22").Value = 5999.3335
$ws.Range("J22").Value = 5999.3335
$ws.Range("L22").Value = 5999.3335
$ws.Range("N22").Value = -6585.3335

$ws.Range("H24").Value = 5000000
$ws.Range("I24").Value = 5000000
$ws.Range("K24").Value = 5000000
$ws.Range("M24").Value = -4999770

$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = $null

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = $null

$ws.Range("H33").Value = 16999
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = $null

$ws.Range("H36").Value = 16999
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = $null

$ws.Range("H37").Value = 30000
$ws.Range("I37").Value = 30000
$ws.Range("K37").Value = 30000
$ws.Range("M37").Value = -29797

$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 10000
$ws.Range("K40").Value = 10000
$ws.Range("M40").Value = -9851

$ws.Range("H52").Value = 10024500
$ws.Range("I52").Value = 20000000
$ws.Range("J52").Value = 49000
$ws.Range("K52").Value = 20000000
$ws.Range("L52").Value = 49000
$ws.Range("M52").Value = -19999774
$ws.Range("N52").Value = -49452

$ws.Range("H58").Value = 20000000
$ws.Range("I58").Value = 20000000
$ws.Range("K58").Value = 20000000
$ws.Range("M58").Value = -19999692

$ws.Range("H62").Value = 4106.25
$ws.Range("I62").Value = 3837.5
$ws.Range("J62").Value = 4375
$ws.Range("K62").Value = 3837.5
$ws.Range("L62").Value = 4375
$ws.Range("M62").Value = -3213.5
$ws.Range("N62").Value = -5623

$ws.Range("H65").Value = 4106.25
$ws.Range("I65").Value = 3837.5
$ws.Range("J65").Value = 4375
$ws.Range("K65").Value = 19187.5
$ws.Range("L65").Value = 21875
$ws.Range("M65").Value = -16067.5
$ws.Range("N65").Value = -28115

$ws.Range("H136").Value = 1100
$ws.Range("I136").Value = 900
$ws.Range("K136").Value = 2700
$ws.Range("M136").Value = -150
